$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 358.3846
$ws.Range("I4").Value = 245.7
$ws.Range("J4").Value = 734
$ws.Range("K4").Value = 245.7
$ws.Range("L4").Value = 734
$ws.Range("M4").Value = -131.7
$ws.Range("N4").Value = -962

$ws.Range("H12").Value = 406.66666
$ws.Range("I12").Value = 110
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 110
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = -1340

$ws.Range("H40").Value = 1979.25
$ws.Range("I40").Value = 1984.7192
$ws.Range("J40").Value = 1950.909
$ws.Range("K40").Value = 1984.7192
$ws.Range("L40").Value = 1950.909
$ws.Range("M40").Value = -1809.7192
$ws.Range("N40").Value = -2300.909

$ws.Range("H70").Value = 2251.5264
$ws.Range("I70").Value = 2322.4443
$ws.Range("J70").Value = 2187.7
$ws.Range("K70").Value = 6967.3329
$ws.Range("L70").Value = 6563.099999999999
$ws.Range("M70").Value = -6697.3329
$ws.Range("N70").Value = -7103.099999999999

$ws.Range("H73").Value = 2251.5264
$ws.Range("I73").Value = 2322.4443
$ws.Range("J73").Value = 2187.7
$ws.Range("K73").Value = 6967.3329
$ws.Range("L73").Value = 6563.099999999999
$ws.Range("M73").Value = -6031.3329
$ws.Range("N73").Value = -8435.099999999999

$ws.Range("H82").Value = 1022.4
$ws.Range("I82").Value = 1022.4
$ws.Range("K82").Value = 3067.2
$ws.Range("M82").Value = -2661.2

$ws.Range("H85").Value = 1022.4
$ws.Range("I85").Value = 1022.4
$ws.Range("K85").Value = 3067.2
$ws.Range("M85").Value = -1663.2

$ws.Range("H138").Value = 1730.88
$ws.Range("I138").Value = 1313.5264
$ws.Range("J138").Value = 3052.5
$ws.Range("K138").Value = 3940.5792
$ws.Range("L138").Value = 9157.5
$ws.Range("M138").Value = 1199.4208
$ws.Range("N138").Value = -19437.5

$ws.Range("H141").Value = 2456.875
$ws.Range("I141").Value = 1859.1666
$ws.Range("J141").Value = 4250
$ws.Range("K141").Value = 5577.4998
$ws.Range("L141").Value = 12750
$ws.Range("M141").Value = -397.4997999999996
$ws.Range("N141").Value = -23110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 149
$ws.Range("I5").Value = 150
$ws.Range("J5").Value = 148
$ws.Range("K5").Value = 150
$ws.Range("L5").Value = 148
$ws.Range("M5").Value = -38
$ws.Range("N5").Value = -372

$ws.Range("H26").Value = 3558
$ws.Range("I26").Value = 596.6667
$ws.Range("J26").Value = 8000
$ws.Range("K26").Value = 596.6667
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = -266.6667
$ws.Range("N26").Value = -8660

$ws.Range("H39").Value = 7475
$ws.Range("I39").Value = 5950
$ws.Range("K39").Value = 5950
$ws.Range("M39").Value = -5430

$ws.Range("H97").Value = 1052.2084
$ws.Range("I97").Value = 740.5
$ws.Range("J97").Value = 1675.625
$ws.Range("K97").Value = 740.5
$ws.Range("L97").Value = 1675.625
$ws.Range("M97").Value = -244.5
$ws.Range("N97").Value = -2667.625

$ws.Range("H132").Value = 3350.1936
$ws.Range("I132").Value = 1835.125
$ws.Range("J132").Value = 4966.2666
$ws.Range("K132").Value = 5505.375
$ws.Range("L132").Value = 14898.7998
$ws.Range("M132").Value = -2975.375
$ws.Range("N132").Value = -19958.7998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 149
$ws.Range("I4").Value = 150
$ws.Range("J4").Value = 148
$ws.Range("K4").Value = 150
$ws.Range("L4").Value = 148
$ws.Range("M4").Value = -35
$ws.Range("N4").Value = -378

$ws.Range("H20").Value = 1648.9333
$ws.Range("I20").Value = 1525
$ws.Range("J20").Value = 1710.9
$ws.Range("K20").Value = 1525
$ws.Range("L20").Value = 1710.9
$ws.Range("M20").Value = -1278
$ws.Range("N20").Value = -2204.9

$ws.Range("H22").Value = 552.6
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 565.75
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 565.75
$ws.Range("M22").Value = -327
$ws.Range("N22").Value = -911.75

$ws.Range("H86").Value = 1519.2307
$ws.Range("I86").Value = 1505.2354
$ws.Range("J86").Value = 1545.6666
$ws.Range("K86").Value = 1505.2354
$ws.Range("L86").Value = 1545.6666
$ws.Range("M86").Value = -382.2354
$ws.Range("N86").Value = -3791.6666

$ws.Range("H89").Value = 1519.2307
$ws.Range("I89").Value = 1505.2354
$ws.Range("J89").Value = 1545.6666
$ws.Range("K89").Value = 7526.177
$ws.Range("L89").Value = 7728.333000000001
$ws.Range("M89").Value = -1910.177
$ws.Range("N89").Value = -18960.333

$ws.Range("H94").Value = 1051.5526
$ws.Range("I94").Value = 563.3333
$ws.Range("J94").Value = 2249.9092
$ws.Range("K94").Value = 563.3333
$ws.Range("L94").Value = 2249.9092
$ws.Range("M94").Value = -112.3333
$ws.Range("N94").Value = -3151.9092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 30643
$ws.Range("J28").Value = 30643
$ws.Range("L28").Value = 30643
$ws.Range("N28").Value = -31133

$ws.Range("H35").Value = 1263.375
$ws.Range("I35").Value = 1263.375
$ws.Range("K35").Value = 1263.375
$ws.Range("M35").Value = -969.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 6711.095
$ws.Range("J124").Value = 6896.65
$ws.Range("L124").Value = 20689.95
$ws.Range("N124").Value = -30509.95

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5620.1113
$ws.Range("I70").Value = 5683.421
$ws.Range("J70").Value = 5469.75
$ws.Range("K70").Value = 5683.421
$ws.Range("L70").Value = 5469.75
$ws.Range("M70").Value = -5413.421
$ws.Range("N70").Value = -6009.75

$ws.Range("H73").Value = 5620.1113
$ws.Range("I73").Value = 5683.421
$ws.Range("J73").Value = 5469.75
$ws.Range("K73").Value = 5683.421
$ws.Range("L73").Value = 5469.75
$ws.Range("M73").Value = -4747.421
$ws.Range("N73").Value = -7341.75

$ws.Range("H80").Value = 3557
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 3779.8
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3779.8
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -5775.8

$ws.Range("H83").Value = 3557
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 3779.8
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 18899
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -28883

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 641.8
$ws.Range("I9").Value = 99
$ws.Range("J9").Value = 1003.6667
$ws.Range("K9").Value = 99
$ws.Range("L9").Value = 1003.6667
$ws.Range("M9").Value = 125
$ws.Range("N9").Value = -1451.6667

$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 1500
$ws.Range("K16").Value = 1500
$ws.Range("M16").Value = -1330

$ws.Range("H46").Value = 1067.1428
$ws.Range("I46").Value = 834.2857
$ws.Range("J46").Value = 1300
$ws.Range("K46").Value = 834.2857
$ws.Range("L46").Value = 1300
$ws.Range("M46").Value = -646.2857
$ws.Range("N46").Value = -1676

$ws.Range("H82").Value = 387550.47
$ws.Range("I82").Value = 835060.8
$ws.Range("J82").Value = 3970.1428
$ws.Range("K82").Value = 835060.8
$ws.Range("L82").Value = 3970.1428
$ws.Range("M82").Value = -834699.8
$ws.Range("N82").Value = -4692.1428

$ws.Range("H85").Value = 387550.47
$ws.Range("I85").Value = 835060.8
$ws.Range("J85").Value = 3970.1428
$ws.Range("K85").Value = 835060.8
$ws.Range("L85").Value = 3970.1428
$ws.Range("M85").Value = -833812.8
$ws.Range("N85").Value = -6466.1428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1114.1333
$ws.Range("I113").Value = 790
$ws.Range("J113").Value = 1600.3334
$ws.Range("K113").Value = 2370
$ws.Range("L113").Value = 4801.0002
$ws.Range("M113").Value = -200
$ws.Range("N113").Value = -9141.0002
